$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.698.55"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").Value = "1.636.54"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'213.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").Value = "'0.503"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.33%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +2.72%  "

$ws.Range("D9").Value = "'0.0623"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("D10").Value = "'19.20"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  +3.35%  "

$ws.Range("D12").Value = "1.865.63"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").Value = "1.631.90"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").Value = "'0.527"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").Value = "26.700.64"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").Value = "'63.27"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("E18").Value = "  +2.49%  "

$ws.Range("D19").Value = "'217.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.06%  "

$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "'4.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "'9.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("E23").Value = "  +2.85%  "

$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").Value = "'148.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.24%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D28").Value = "'6.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.07%  "

$ws.Range("E29").Value = "  +1.85%  "

$ws.Range("D30").Value = "'0.0509"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.04%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").Value = "'3.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.42%  "

$ws.Range("D33").Value = "'2.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.40%  "

$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("D36").Value = "1.202.76"
$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("D37").Value = "'0.0173"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.13%  "

$ws.Range("D38").Value = "'0.812"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").Value = "'0.507"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.37%  "

$ws.Range("D41").Value = "'2.28"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("D42").Value = "'5.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.70%  "

$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D44").Value = "1.773.52"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("D45").Value = "'92.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("D46").Value = "'1.56"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.55%  "

$ws.Range("D47").Value = "'54.84"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("D49").Value = "'7.65"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.28%  "

$ws.Range("D50").Value = "'0.410"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("E51").Value = "  +0.11%  "
